# Fixed naive component forecaster bug - Presentation state 11.02.
#
# Each data row (r = 2..24) holds a rolling window of QoQ forecast errors
# laid out left-to-right from column B onward, newest first. The naive
# forecaster was writing the freshly computed error into the wrong slot;
# the fix re-aligns the row by pushing every existing value one column to
# the right (B->C, C->D, ...) and writing the newly computed error into the
# now-vacated column B. Rows that were already full (out to column K) drop
# their oldest (rightmost) observation as it rolls off the window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly computed naive QoQ error for each row, destined for column B.
$newFirstValue = @{
    2  = 2.297389002388887
    3  = 8.826710628892494
    4  = -9.780318414391347
    5  = -1.200275438764269
    6  = 0.3719860057927588
    7  = -2.702915518772638
    8  = -0.2307826431404359
    9  = -0.5654386276933741
    10 = -0.6603092772102132
    11 = -0.15162438770796
    12 = -0.2053460154962278
    13 = 0.6162032393936197
    14 = 1.652643173475852
    15 = 0.3110387314724781
    16 = 0.2388379152847414
    17 = 0.6508000635779043
    18 = 0.2387740594105157
    19 = 0.3465902496671606
    20 = 0.00230005330798793
    21 = -0.1902738424076751
    22 = -0.3325070745318338
    23 = 0.1656141382254278
    24 = -0.09587373626955231
}

# Column B is 2; the window caps at column K (index 11), i.e. 10 slots.
$firstCol = 2
$lastCol = 11
$maxSlots = $lastCol - $firstCol + 1

foreach ($r in 2..24) {
    # Read the row's existing values (old column B .. old last populated column).
    $existing = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq $null) { break }
        $existing += $cell.Value2
    }

    # Shift the existing values one column to the right, dropping the
    # oldest entry once the row would otherwise overflow past column K.
    $shifted = $existing
    if ($shifted.Length -gt ($maxSlots - 1)) {
        $shifted = $shifted[0..($maxSlots - 2)]
    }

    # Column B gets the newly computed value.
    $ws.Cells.Item($r, $firstCol).Value2 = $newFirstValue[$r]

    # Write the shifted values back starting at column C.
    for ($i = 0; $i -lt $shifted.Length; $i++) {
        $ws.Cells.Item($r, $firstCol + 1 + $i).Value2 = $shifted[$i]
    }
}
